$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price values in column D stay as text, matching the
# original inline-string cell type used throughout this sheet.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "68.018.33"
$ws.Range("E2").Value = "  +0.48%  "
$ws.Range("D3").Value = "3.321.96"
$ws.Range("E3").Value = "  +1.83%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "187.29"
$ws.Range("E5").Value = "  +1.66%  "
$ws.Range("D6").Value = "583.17"
$ws.Range("E6").Value = "  +0.59%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "0.599"
$ws.Range("E8").Value = "  -0.29%  "
$ws.Range("D9").Value = "0.130"
$ws.Range("E9").Value = "  +0.68%  "
$ws.Range("D10").Value = "6.69"
$ws.Range("E10").Value = "  +1.66%  "
$ws.Range("D11").Value = "0.408"
$ws.Range("E11").Value = "  +0.21%  "
$ws.Range("D12").Value = "3.902.23"
$ws.Range("E12").Value = "  +2.04%  "
$ws.Range("E13").Value = "  -1.73%  "
$ws.Range("D14").Value = "27.69"
$ws.Range("E14").Value = "  +1.29%  "
$ws.Range("D15").Value = "68.156.86"
$ws.Range("E15").Value = "  +0.59%  "
$ws.Range("E16").Value = "  +0.31%  "
$ws.Range("D17").Value = "3.315.39"
$ws.Range("E17").Value = "  +1.94%  "
$ws.Range("D18").Value = "446.38"
$ws.Range("E18").Value = "  +13.09%  "
$ws.Range("D19").Value = "5.74"
$ws.Range("E19").Value = "  +0.62%  "
$ws.Range("D20").Value = "13.60"
$ws.Range("E20").Value = "  +1.33%  "
$ws.Range("D21").Value = "7.73"
$ws.Range("E21").Value = "  +2.54%  "
$ws.Range("D22").Value = "75.45"
$ws.Range("E22").Value = "  +6.60%  "
$ws.Range("D23").Value = "0.996"
$ws.Range("E23").Value = "  -0.48%  "
$ws.Range("B24").Value = "WrappedeETH"
$ws.Range("C24").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D24").Value = "3.476.31"
$ws.Range("B25").Value = "Polygon"
$ws.Range("C25").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D25").Value = "0.517"
$ws.Range("E25").Value = "  +1.91%  "
$ws.Range("D26").Value = "0.0000119"
$ws.Range("E26").Value = "  +1.92%  "
$ws.Range("E27").Value = "  +0.83%  "
$ws.Range("D28").Value = "9.22"
$ws.Range("E28").Value = "  -2.60%  "
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.38%  "
$ws.Range("D30").Value = "1.99"
$ws.Range("E30").Value = "  +2.17%  "
$ws.Range("D31").Value = "23.01"
$ws.Range("E31").Value = "  +2.03%  "
$ws.Range("D32").Value = "5.39"
$ws.Range("E32").Value = "  -0.96%  "
$ws.Range("D33").Value = "1.25"
$ws.Range("E33").Value = "  +1.03%  "
$ws.Range("B34").Value = "USDe"
$ws.Range("C34").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  +0.06%  "
$ws.Range("B35").Value = "Aptos"
$ws.Range("C35").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D35").Value = "6.83"
$ws.Range("E35").Value = "  -1.10%  "
$ws.Range("E36").Value = "  +5.32%  "
$ws.Range("D37").Value = "163.82"
$ws.Range("E37").Value = "  -0.17%  "
$ws.Range("D38").Value = "1.89"
$ws.Range("E38").Value = "  +0.48%  "
$ws.Range("D39").Value = "27.05"
$ws.Range("E39").Value = "  +1.12%  "
$ws.Range("D40").Value = "4.55"
$ws.Range("E40").Value = "  +1.39%  "
$ws.Range("D41").Value = "0.791"
$ws.Range("E41").Value = "  -1.74%  "
$ws.Range("D42").Value = "6.41"
$ws.Range("E42").Value = "  +2.25%  "
$ws.Range("D43").Value = "2.697.96"
$ws.Range("E43").Value = "  +1.54%  "
$ws.Range("D44").Value = "40.82"
$ws.Range("E44").Value = "  +0.14%  "
$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").Value = "0.0679"
$ws.Range("E45").Value = "  +0.09%  "
$ws.Range("B46").Value = "dogwifhat"
$ws.Range("C46").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D46").Value = "2.42"
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("D47").Value = "24.67"
$ws.Range("E47").Value = "  +0.37%  "
$ws.Range("D48").Value = "327.79"
$ws.Range("E48").Value = "  -2.02%  "
$ws.Range("D49").Value = "0.0277"
$ws.Range("E49").Value = "  +1.41%  "
$ws.Range("D50").Value = "31.79"
$ws.Range("E50").Value = "  +4.44%  "
$ws.Range("D51").Value = "0.991"
$ws.Range("E51").Value = "  +2.16%  "
